# Regenerate save_data to use K (strikeouts) instead of Strike# for column G
# Updates G2:G9 on the active sheet with the recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 3
    6 = 2
    7 = 1
    8 = 1
    9 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
